# Generate Report for Archive
#
# 1) Update the localization status text "Ready for handoff" -> "In Translation"
#    wherever it appears (Overview sheet columns E/F, and the per-locale
#    "Status" column C on the zh-cn / de-de sheets).
# 2) Narrow the "Status" columns (Overview!E:F and Status column C on the
#    locale sheets) from their old width down to the new, narrower width.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Sheet 1: Overview ---
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# Narrow columns E and F (status columns) on the Overview sheet.
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5

# --- Sheet 2: zh-cn ---
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# Narrow the Status column (C) on the zh-cn sheet.
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5

# --- Sheet 3: de-de ---
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# Narrow the Status column (C) on the de-de sheet.
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
